# Auto-generated edit script: updates profit/price figures across multiple sheets
# per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016278.2
$ws.Range("I38").Value = 2016278.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 6048834.6
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -6048462.6
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 2126.2222
$ws.Range("I43").Value = 10001
$ws.Range("J43").Value = 1141.875
$ws.Range("K43").Value = 10001
$ws.Range("L43").Value = 1141.875
$ws.Range("M43").Value = -9932
$ws.Range("N43").Value = -1279.875
$ws.Range("H112").Value = 1086.6471
$ws.Range("J112").Value = 1165.5333
$ws.Range("L112").Value = 3496.5999
$ws.Range("N112").Value = -5712.5999
$ws.Range("H116").Value = 3660
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3660
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3660
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -10544

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6750
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = -2885
$ws.Range("N3").Value = -8230
$ws.Range("H32").Value = 5726.628
$ws.Range("I32").Value = 5160.4287
$ws.Range("J32").Value = 29507
$ws.Range("K32").Value = 5160.4287
$ws.Range("L32").Value = 29507
$ws.Range("M32").Value = -4873.4287
$ws.Range("N32").Value = -30081
$ws.Range("H45").Value = 84827.164
$ws.Range("I45").Value = 126063.125
$ws.Range("J45").Value = 2355.25
$ws.Range("K45").Value = 126063.125
$ws.Range("L45").Value = 2355.25
$ws.Range("M45").Value = -125686.125
$ws.Range("N45").Value = -3109.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 34095.902
$ws.Range("I20").Value = 45535.523
$ws.Range("J20").Value = 1207
$ws.Range("K20").Value = 45535.523
$ws.Range("L20").Value = 1207
$ws.Range("M20").Value = -45288.523
$ws.Range("N20").Value = -1701
$ws.Range("H86").Value = 139294.5
$ws.Range("I86").Value = 158972.28
$ws.Range("J86").Value = 1550
$ws.Range("K86").Value = 158972.28
$ws.Range("L86").Value = 1550
$ws.Range("M86").Value = -157849.28
$ws.Range("N86").Value = -3796
$ws.Range("H89").Value = 139294.5
$ws.Range("I89").Value = 158972.28
$ws.Range("J89").Value = 1550
$ws.Range("K89").Value = 794861.4
$ws.Range("L89").Value = 7750
$ws.Range("M89").Value = -789245.4
$ws.Range("N89").Value = -18982
$ws.Range("H105").Value = 126931.375
$ws.Range("I105").Value = 251943.5
$ws.Range("J105").Value = 1919.25
$ws.Range("K105").Value = 251943.5
$ws.Range("L105").Value = 1919.25
$ws.Range("M105").Value = -250196.5
$ws.Range("N105").Value = -5413.25
$ws.Range("H134").Value = 3274.44
$ws.Range("I134").Value = 3024.8333
$ws.Range("J134").Value = 3916.2856
$ws.Range("K134").Value = 9074.499899999999
$ws.Range("L134").Value = 11748.8568
$ws.Range("M134").Value = -6539.499899999999
$ws.Range("N134").Value = -16818.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1329.6666
$ws.Range("I16").Value = 1329.6666
$ws.Range("K16").Value = 1329.6666
$ws.Range("M16").Value = -1042.6666
$ws.Range("H31").Value = 2879.7646
$ws.Range("I31").Value = 1783.25
$ws.Range("J31").Value = 3587.1936
$ws.Range("K31").Value = 1783.25
$ws.Range("L31").Value = 3587.1936
$ws.Range("M31").Value = -1488.25
$ws.Range("N31").Value = -4177.193600000001
$ws.Range("H34").Value = 2879.7646
$ws.Range("I34").Value = 1783.25
$ws.Range("J34").Value = 3587.1936
$ws.Range("K34").Value = 1783.25
$ws.Range("L34").Value = 3587.1936
$ws.Range("M34").Value = -1581.25
$ws.Range("N34").Value = -3991.1936
$ws.Range("H45").Value = 14666
$ws.Range("I45").Value = 14499
$ws.Range("K45").Value = 14499
$ws.Range("M45").Value = -13906
$ws.Range("H74").Value = 27081.555
$ws.Range("J74").Value = 27081.555
$ws.Range("L74").Value = 27081.555
$ws.Range("N74").Value = -28829.555
$ws.Range("H77").Value = 27081.555
$ws.Range("J77").Value = 27081.555
$ws.Range("L77").Value = 81244.66500000001
$ws.Range("N77").Value = -89980.66500000001
$ws.Range("H86").Value = 4352
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4352
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4352
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -6598
$ws.Range("H88").Value = 24007.25
$ws.Range("J88").Value = 24007.25
$ws.Range("L88").Value = 24007.25
$ws.Range("N88").Value = -24819.25
$ws.Range("H89").Value = 4352
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4352
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 21760
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -32992
$ws.Range("H91").Value = 24007.25
$ws.Range("J91").Value = 24007.25
$ws.Range("L91").Value = 24007.25
$ws.Range("N91").Value = -26815.25
$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992
$ws.Range("H96").Value = 3508
$ws.Range("J96").Value = 3508
$ws.Range("L96").Value = 3508
$ws.Range("N96").Value = -9000
$ws.Range("H113").Value = 1329.6666
$ws.Range("I113").Value = 1329.6666
$ws.Range("K113").Value = 1329.6666
$ws.Range("M113").Value = 840.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 700
$ws.Range("J34").Value = 700
$ws.Range("L34").Value = 2100
$ws.Range("N34").Value = -2268
$ws.Range("H131").Value = 776.55554
$ws.Range("J131").Value = 817.47253
$ws.Range("L131").Value = 2452.41759
$ws.Range("N131").Value = -12532.41759
$ws.Range("H132").Value = 1533.225
$ws.Range("I132").Value = 701.1875
$ws.Range("J132").Value = 2087.9167
$ws.Range("K132").Value = 6310.6875
$ws.Range("L132").Value = 18791.2503
$ws.Range("M132").Value = -3780.6875
$ws.Range("N132").Value = -23851.2503
$ws.Range("H139").Value = 2027.8334
$ws.Range("J139").Value = 3375.8
$ws.Range("L139").Value = 10127.4
$ws.Range("N139").Value = -20407.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3748
$ws.Range("H132").Value = 3652.9546
$ws.Range("I132").Value = 3334.625
$ws.Range("K132").Value = 10003.875
$ws.Range("M132").Value = -7473.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 16400
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 16400
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 16400
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -16786

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1277.909
$ws.Range("I136").Value = 490.7742
$ws.Range("J136").Value = 2294.625
$ws.Range("K136").Value = 1472.3226
$ws.Range("L136").Value = 6883.875
$ws.Range("M136").Value = 1077.6774
$ws.Range("N136").Value = -11983.875

Write-Output "Applied scheduled data refresh to Aegis_Profits sheets."
